# Refresh cached Market Board valuation columns (H-N) across every Disciple of the
# Hand/Land profession sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the latest
# scheduled-runner price pull. Leve/recipe reference columns (A-G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1041.1052
$ws.Range("I98").Value = 990.25714
$ws.Range("K98").Value = 990.25714
$ws.Range("M98").Value = 507.74286
$ws.Range("H112").Value = 2625.9546
$ws.Range("J112").Value = 2625.9546
$ws.Range("L112").Value = 7877.8638
$ws.Range("N112").Value = -10093.8638
$ws.Range("H122").Value = 1041.1052
$ws.Range("I122").Value = 990.25714
$ws.Range("K122").Value = 2970.77142
$ws.Range("M122").Value = -520.77142
$ws.Range("H132").Value = 1283.826
$ws.Range("I132").Value = 1166.6
$ws.Range("J132").Value = 2065.3333
$ws.Range("K132").Value = 3499.8
$ws.Range("L132").Value = 6195.999899999999
$ws.Range("M132").Value = -969.7999999999997
$ws.Range("N132").Value = -11255.9999
$ws.Range("H137").Value = 38594.48
$ws.Range("I137").Value = 1384.1875
$ws.Range("J137").Value = 92718.55
$ws.Range("K137").Value = 4152.5625
$ws.Range("L137").Value = 278155.65
$ws.Range("M137").Value = -1602.5625
$ws.Range("N137").Value = -283255.65
$ws.Range("H138").Value = 3585.4187
$ws.Range("I138").Value = 3712.3635
$ws.Range("J138").Value = 3452.4285
$ws.Range("K138").Value = 11137.0905
$ws.Range("L138").Value = 10357.2855
$ws.Range("M138").Value = -5997.0905
$ws.Range("N138").Value = -20637.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4709
$ws.Range("I31").Value = 4709
$ws.Range("K31").Value = 4709
$ws.Range("M31").Value = -4415
$ws.Range("H32").Value = 14553.981
$ws.Range("I32").Value = 9644.629999999999
$ws.Range("J32").Value = 19463.334
$ws.Range("K32").Value = 9644.629999999999
$ws.Range("L32").Value = 19463.334
$ws.Range("M32").Value = -9357.629999999999
$ws.Range("N32").Value = -20037.334
$ws.Range("H61").Value = 22237.365
$ws.Range("I61").Value = 38462
$ws.Range("K61").Value = 38462
$ws.Range("M61").Value = -38250
$ws.Range("H74").Value = 813.35
$ws.Range("I74").Value = 570.5
$ws.Range("J74").Value = 2999
$ws.Range("K74").Value = 570.5
$ws.Range("L74").Value = 2999
$ws.Range("M74").Value = 303.5
$ws.Range("N74").Value = -4747
$ws.Range("H77").Value = 813.35
$ws.Range("I77").Value = 570.5
$ws.Range("J77").Value = 2999
$ws.Range("K77").Value = 2852.5
$ws.Range("L77").Value = 14995
$ws.Range("M77").Value = 1515.5
$ws.Range("N77").Value = -23731
$ws.Range("H88").Value = 3827.818
$ws.Range("I88").Value = 3001.5
$ws.Range("K88").Value = 3001.5
$ws.Range("M88").Value = -2595.5
$ws.Range("H91").Value = 3827.818
$ws.Range("I91").Value = 3001.5
$ws.Range("K91").Value = 3001.5
$ws.Range("M91").Value = -1597.5
$ws.Range("H122").Value = 23946.295
$ws.Range("I122").Value = 23946.295
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 71838.88499999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -69388.88499999999
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 22237.365
$ws.Range("I136").Value = 38462
$ws.Range("K136").Value = 115386
$ws.Range("M136").Value = -112836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 224107.78
$ws.Range("I86").Value = 1852.8572
$ws.Range("J86").Value = 1002000
$ws.Range("K86").Value = 1852.8572
$ws.Range("L86").Value = 1002000
$ws.Range("M86").Value = -729.8571999999999
$ws.Range("N86").Value = -1004246
$ws.Range("H89").Value = 224107.78
$ws.Range("I89").Value = 1852.8572
$ws.Range("J89").Value = 1002000
$ws.Range("K89").Value = 9264.286
$ws.Range("L89").Value = 5010000
$ws.Range("M89").Value = -3648.286
$ws.Range("N89").Value = -5021232
$ws.Range("H99").Value = 635.375
$ws.Range("I99").Value = 657.5714
$ws.Range("K99").Value = 657.5714
$ws.Range("M99").Value = 840.4286
$ws.Range("H107").Value = 1915.3334
$ws.Range("I107").Value = 1898.4
$ws.Range("K107").Value = 1898.4
$ws.Range("M107").Value = 21.59999999999991
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2787.25
$ws.Range("I31").Value = 1761.375
$ws.Range("J31").Value = 4839
$ws.Range("K31").Value = 1761.375
$ws.Range("L31").Value = 4839
$ws.Range("M31").Value = -1466.375
$ws.Range("N31").Value = -5429
$ws.Range("H34").Value = 2787.25
$ws.Range("I34").Value = 1761.375
$ws.Range("J34").Value = 4839
$ws.Range("K34").Value = 1761.375
$ws.Range("L34").Value = 4839
$ws.Range("M34").Value = -1559.375
$ws.Range("N34").Value = -5243
$ws.Range("H105").Value = 640.75
$ws.Range("I105").Value = 581.6667
$ws.Range("J105").Value = 818
$ws.Range("K105").Value = 581.6667
$ws.Range("L105").Value = 818
$ws.Range("M105").Value = 1165.3333
$ws.Range("N105").Value = -4312
$ws.Range("H107").Value = 1379.6072
$ws.Range("I107").Value = 1332.1305
$ws.Range("K107").Value = 1332.1305
$ws.Range("M107").Value = 587.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125.958336
$ws.Range("I2").Value = 192.33333
$ws.Range("K2").Value = 1153.99998
$ws.Range("M2").Value = -1040.99998
$ws.Range("H80").Value = 2422.625
$ws.Range("J80").Value = 2422.625
$ws.Range("L80").Value = 7267.875
$ws.Range("N80").Value = -9139.875
$ws.Range("H83").Value = 2422.625
$ws.Range("J83").Value = 2422.625
$ws.Range("L83").Value = 21803.625
$ws.Range("N83").Value = -31163.625
$ws.Range("H101").Value = 6166.6665
$ws.Range("J101").Value = 6166.6665
$ws.Range("L101").Value = 18499.9995
$ws.Range("N101").Value = -23367.9995
$ws.Range("H107").Value = 734
$ws.Range("I107").Value = 691.5
$ws.Range("J107").Value = 768
$ws.Range("K107").Value = 2074.5
$ws.Range("L107").Value = 2304
$ws.Range("M107").Value = -154.5
$ws.Range("N107").Value = -6144
$ws.Range("H127").Value = 1644
$ws.Range("J127").Value = 1644
$ws.Range("L127").Value = 4932
$ws.Range("N127").Value = -14852
$ws.Range("H131").Value = 15139.692
$ws.Range("J131").Value = 16642.236
$ws.Range("L131").Value = 49926.708
$ws.Range("N131").Value = -60006.708
$ws.Range("H137").Value = 5703.643
$ws.Range("I137").Value = 1370.5
$ws.Range("K137").Value = 4111.5
$ws.Range("M137").Value = 988.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2929.923
$ws.Range("I80").Value = 2999
$ws.Range("K80").Value = 2999
$ws.Range("M80").Value = -2001
$ws.Range("H83").Value = 2929.923
$ws.Range("I83").Value = 2999
$ws.Range("K83").Value = 14995
$ws.Range("M83").Value = -10003
$ws.Range("H113").Value = 1551.4286
$ws.Range("I113").Value = 1165.25
$ws.Range("J113").Value = 2066.3333
$ws.Range("K113").Value = 1165.25
$ws.Range("L113").Value = 2066.3333
$ws.Range("M113").Value = 1004.75
$ws.Range("N113").Value = -6406.3333
$ws.Range("H126").Value = 2097461.2
$ws.Range("I126").Value = 2648440.8
$ws.Range("K126").Value = 7945322.399999999
$ws.Range("M126").Value = -7942852.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2995.158
$ws.Range("J7").Value = 3045.9092
$ws.Range("L7").Value = 3045.9092
$ws.Range("N7").Value = -3269.9092
$ws.Range("H40").Value = 9977.041999999999
$ws.Range("I40").Value = 8552.777
$ws.Range("J40").Value = 14249.833
$ws.Range("K40").Value = 8552.777
$ws.Range("L40").Value = 14249.833
$ws.Range("M40").Value = -8416.777
$ws.Range("N40").Value = -14521.833
$ws.Range("H126").Value = 2995.158
$ws.Range("J126").Value = 3045.9092
$ws.Range("L126").Value = 9137.7276
$ws.Range("N126").Value = -14077.7276
$ws.Range("H132").Value = 4181.706
$ws.Range("I132").Value = 4621
$ws.Range("K132").Value = 13863
$ws.Range("M132").Value = -11333
$ws.Range("H136").Value = 1436.091
$ws.Range("I136").Value = 1436.091
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4308.272999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1758.272999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 22222
$ws.Range("J29").Value = 22222
$ws.Range("L29").Value = 22222
$ws.Range("N29").Value = -22802
$ws.Range("H81").Value = 749.75
$ws.Range("I81").Value = 749.75
$ws.Range("K81").Value = 1499.5
$ws.Range("M81").Value = -438.5
$ws.Range("H84").Value = 749.75
$ws.Range("I84").Value = 749.75
$ws.Range("K84").Value = 7497.5
$ws.Range("M84").Value = -2193.5
$ws.Range("H113").Value = 988.7692
$ws.Range("I113").Value = 887.1667
$ws.Range("J113").Value = 1075.8572
$ws.Range("K113").Value = 2661.5001
$ws.Range("L113").Value = 3227.5716
$ws.Range("M113").Value = -491.5001000000002
$ws.Range("N113").Value = -7567.571599999999
$ws.Range("H122").Value = 33308.117
$ws.Range("I122").Value = 39572
$ws.Range("J122").Value = 6999.8
$ws.Range("K122").Value = 118716
$ws.Range("L122").Value = 20999.4
$ws.Range("M122").Value = -116266
$ws.Range("N122").Value = -25899.4
$ws.Range("H123").Value = 44469.555
$ws.Range("J123").Value = 44469.555
$ws.Range("L123").Value = 44469.555
$ws.Range("N123").Value = -54269.555
$ws.Range("H136").Value = 12347257
$ws.Range("I136").Value = 17922370
$ws.Range("K136").Value = 53767110
$ws.Range("M136").Value = -53764560

